# Auto-generated script to update Kraken market price cache values
# across all sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3666.3333
$ws.Range("I18").Value = 2749.5
$ws.Range("K18").Value = 2749.5
$ws.Range("M18").Value = -2465.5
$ws.Range("H43").Value = 1498.75
$ws.Range("I43").Value = 1995
$ws.Range("J43").Value = 1333.3334
$ws.Range("K43").Value = 1995
$ws.Range("L43").Value = 1333.3334
$ws.Range("M43").Value = -1926
$ws.Range("N43").Value = -1471.3334
$ws.Range("H97").Value = 3762.5
$ws.Range("J97").Value = 3762.5
$ws.Range("L97").Value = 11287.5
$ws.Range("N97").Value = -12279.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 50000
$ws.Range("I7").Value = 50000
$ws.Range("K7").Value = 50000
$ws.Range("M7").Value = -49886
$ws.Range("H32").Value = 4361.5386
$ws.Range("I32").Value = 4361.5386
$ws.Range("K32").Value = 4361.5386
$ws.Range("M32").Value = -4074.5386
$ws.Range("H122").Value = 7880.0586
$ws.Range("I122").Value = 8613.462
$ws.Range("J122").Value = 5496.5
$ws.Range("K122").Value = 25840.386
$ws.Range("L122").Value = 16489.5
$ws.Range("M122").Value = -23390.386
$ws.Range("N122").Value = -21389.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1030.3334
$ws.Range("I8").Value = 700
$ws.Range("J8").Value = 1360.6666
$ws.Range("K8").Value = 700
$ws.Range("L8").Value = 1360.6666
$ws.Range("M8").Value = -560
$ws.Range("N8").Value = -1640.6666
$ws.Range("H29").Value = 1000
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 1000
$ws.Range("N29").Value = -1578
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -466
$ws.Range("H107").Value = 18185.285
$ws.Range("I107").Value = 7049.6665
$ws.Range("K107").Value = 7049.6665
$ws.Range("M107").Value = -5129.6665
$ws.Range("H134").Value = 8961.166999999999
$ws.Range("I134").Value = 3389
$ws.Range("K134").Value = 10167
$ws.Range("M134").Value = -7632

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 870.4286
$ws.Range("J16").Value = 820.6
$ws.Range("L16").Value = 820.6
$ws.Range("N16").Value = -1394.6
$ws.Range("H31").Value = 5556
$ws.Range("J31").Value = 8398.6
$ws.Range("L31").Value = 8398.6
$ws.Range("N31").Value = -8988.6
$ws.Range("H34").Value = 5556
$ws.Range("J34").Value = 8398.6
$ws.Range("L34").Value = 8398.6
$ws.Range("N34").Value = -8802.6
$ws.Range("H50").Value = 25264
$ws.Range("I50").Value = 26666.666
$ws.Range("J50").Value = 24881.455
$ws.Range("K50").Value = 26666.666
$ws.Range("L50").Value = 24881.455
$ws.Range("M50").Value = -26041.666
$ws.Range("N50").Value = -26131.455
$ws.Range("H51").Value = 21999.4
$ws.Range("J51").Value = 21999.4
$ws.Range("L51").Value = 21999.4
$ws.Range("N51").Value = -23471.4
$ws.Range("H58").Value = 619
$ws.Range("I58").Value = 619
$ws.Range("K58").Value = 619
$ws.Range("M58").Value = -416
$ws.Range("H60").Value = 16931.25
$ws.Range("I60").Value = 9487.5
$ws.Range("K60").Value = 9487.5
$ws.Range("M60").Value = -8976.5
$ws.Range("H61").Value = 21999.4
$ws.Range("J61").Value = 21999.4
$ws.Range("L61").Value = 21999.4
$ws.Range("N61").Value = -22695.4
$ws.Range("H94").Value = 1732.6666
$ws.Range("I94").Value = 1699
$ws.Range("K94").Value = 1699
$ws.Range("M94").Value = -1248
$ws.Range("H107").Value = 816.63635
$ws.Range("I107").Value = 861.8570999999999
$ws.Range("K107").Value = 861.8570999999999
$ws.Range("M107").Value = 1058.1429
$ws.Range("H113").Value = 870.4286
$ws.Range("J113").Value = 820.6
$ws.Range("L113").Value = 820.6
$ws.Range("N113").Value = -5160.6
$ws.Range("H134").Value = 5749.25
$ws.Range("I134").Value = 5000
$ws.Range("J134").Value = 5999
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 17997
$ws.Range("M134").Value = -12465
$ws.Range("N134").Value = -23067
$ws.Range("H136").Value = 619
$ws.Range("I136").Value = 619
$ws.Range("K136").Value = 1857
$ws.Range("M136").Value = 693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 750
$ws.Range("J32").Value = 750
$ws.Range("L32").Value = 2250
$ws.Range("N32").Value = -2816
$ws.Range("H122").Value = 2488.2
$ws.Range("I122").Value = 1998
$ws.Range("J122").Value = 2610.75
$ws.Range("K122").Value = 17982
$ws.Range("L122").Value = 23496.75
$ws.Range("M122").Value = -15532
$ws.Range("N122").Value = -28396.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 120
$ws.Range("I2").Value = 153.75
$ws.Range("J2").Value = 52.5
$ws.Range("K2").Value = 153.75
$ws.Range("L2").Value = 52.5
$ws.Range("M2").Value = -40.75
$ws.Range("N2").Value = -278.5
$ws.Range("H57").Value = 17500
$ws.Range("H102").Value = 1750
$ws.Range("I102").Value = 1750
$ws.Range("K102").Value = 1750
$ws.Range("M102").Value = -128
$ws.Range("H107").Value = 370.6
$ws.Range("I107").Value = 200
$ws.Range("K107").Value = 200
$ws.Range("M107").Value = 1720
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2587.6667
$ws.Range("I122").Value = 2661.125
$ws.Range("K122").Value = 7983.375
$ws.Range("M122").Value = -5533.375
$ws.Range("H132").Value = 3127.6924
$ws.Range("I132").Value = 2420.625
$ws.Range("J132").Value = 4259
$ws.Range("K132").Value = 7261.875
$ws.Range("L132").Value = 12777
$ws.Range("M132").Value = -4731.875
$ws.Range("N132").Value = -17837

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6987.909
$ws.Range("I40").Value = 6985.222
$ws.Range("K40").Value = 6985.222
$ws.Range("M40").Value = -6849.222
$ws.Range("H55").Value = 2601.9285
$ws.Range("I55").Value = 279.5
$ws.Range("J55").Value = 4343.75
$ws.Range("K55").Value = 279.5
$ws.Range("L55").Value = 4343.75
$ws.Range("M55").Value = -106.5
$ws.Range("N55").Value = -4689.75
$ws.Range("H61").Value = 3234.1428
$ws.Range("I61").Value = 3406.5
$ws.Range("K61").Value = 3406.5
$ws.Range("M61").Value = -3204.5
$ws.Range("H113").Value = 3234.1428
$ws.Range("I113").Value = 3406.5
$ws.Range("K113").Value = 3406.5
$ws.Range("M113").Value = -1236.5
$ws.Range("H122").Value = 3290.9167
$ws.Range("I122").Value = 3381.4546
$ws.Range("J122").Value = 2295
$ws.Range("K122").Value = 10144.3638
$ws.Range("L122").Value = 6885
$ws.Range("M122").Value = -7694.363799999999
$ws.Range("N122").Value = -11785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H126").Value = 1915.8334
$ws.Range("I126").Value = 1979.4
$ws.Range("J126").Value = 1598
$ws.Range("K126").Value = 5938.200000000001
$ws.Range("L126").Value = 4794
$ws.Range("M126").Value = -3468.200000000001
$ws.Range("N126").Value = -3474.25
